$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data rows (row 1 header stays the same)
$ws.Range("A2").Value = "mngr168479"
$ws.Range("B2").Value = "vupemYz"

$ws.Range("A3").Value = "mngr83460"
$ws.Range("B3").Value = "qAbUzyj"

$ws.Range("A4").Value = "mngr168479"
$ws.Range("B4").Value = "vupemYz"

$ws.Range("A5").Value = "mngr164225"
$ws.Range("B5").Value = "jahetAp"

$ws.Range("A6").Value = "mngr168479"
$ws.Range("B6").Value = "vupemYz"

# Update the selected cell in the sheet view
$ws.Range("A6").Select()
